$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.30138635635376
$ws.Range("B1").Value = 2.202964305877686
$ws.Range("C1").Value = 5.033305644989014
$ws.Range("D1").Value = 1.983607172966003
$ws.Range("E1").Value = 1.074829816818237
